$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4848494546311366
$ws.Range("C2").Value = 0.9857587014232243

$ws.Range("B3").Value = 0.09668751914543174
$ws.Range("C3").Value = 0.9986637523058749

$ws.Range("B4").Value = 0.03224830283962941
$ws.Range("C4").Value = 0.9996676236188253

$ws.Range("B5").Value = 0.09715399643991034
$ws.Range("C5").Value = 0.9994212813634508
